$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.795.12"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "1.915.90"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'250.52"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'46.44"
$ws.Range("E8").Value = "  +6.89%  "
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("D10").Value = "'58.54"
$ws.Range("E10").Value = "  +9.84%  "
$ws.Range("D11").Value = "'0.0765"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "'0.0999"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "'14.72"
$ws.Range("E13").Value = "  +8.98%  "
$ws.Range("D14").Value = "'0.816"
$ws.Range("E14").Value = "  +7.30%  "
$ws.Range("D15").Value = "2.191.73"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "'5.15"
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("D17").Value = "1.911.01"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "36.755.97"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").Value = "'75.00"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E20").Value = "  +4.47%  "
$ws.Range("D21").Value = "'251.26"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").Value = "'13.42"
$ws.Range("E22").Value = "  +4.90%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'2.64"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'168.19"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "'8.81"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").Value = "'18.78"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'4.57"
$ws.Range("E31").Value = "  +6.13%  "
$ws.Range("D32").Value = "'0.0618"
$ws.Range("D33").Value = "'4.35"
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").Value = "'0.0898"
$ws.Range("E34").Value = "  +23.48%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "'1.52"
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("D38").Value = "'0.876"
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("D39").Value = "'17.92"
$ws.Range("E39").Value = "  +52.19%  "
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("D41").Value = "'105.37"
$ws.Range("E41").Value = "  +8.54%  "
$ws.Range("D42").Value = "'0.0228"
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("D43").Value = "'17.71"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("E44").Value = "  +21.95%  "
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").Value = "1.350.57"
$ws.Range("D47").Value = "'2.39"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "'0.0813"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").Value = "'2.81"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "'6.47"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").Value = "'43.53"
$ws.Range("E51").Value = "  +2.73%  "
